# "Fix bugs, added review, results"
#
# The true_false questions (rows 2, 4, 6, 8, 10, 12, 14, 16) stored their
# "correct" answer as the literal text "true", while every multiple_choice
# question records its correct answer using the "option_N" scheme (matching
# the option_1..option_4 headers). This fixes the true_false rows to use the
# same "option_1" convention so downstream grading code only has to deal with
# one answer-key format.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H2").Value  = "option_1"
$ws.Range("H4").Value  = "option_1"
$ws.Range("H6").Value  = "option_1"
$ws.Range("H8").Value  = "option_1"
$ws.Range("H10").Value = "option_1"
$ws.Range("H12").Value = "option_1"
$ws.Range("H14").Value = "option_1"
$ws.Range("H16").Value = "option_1"

# Leave behind the selection where the author's review ended up.
$ws.Range("K28").Select()
